$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9693716918425304
$ws.Range("J2").Value = 0.9693716918425304
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.502639
$ws.Range("N2").Value = 31.507917
$ws.Range("O2").Value = 0.9701256668284471
$ws.Range("P2").Value = 0.970125666828447
$ws.Range("Q2").Value = 96.17223821568068
$ws.Range("R2").Value = 865.5501439411261
$ws.Range("S2").Value = 0.9404123589533547
$ws.Range("T2").Value = 0.9404123589533546

# Row 3
$ws.Range("I3").Value = 0.9693716918425304
$ws.Range("J3").Value = 0.9693716918425304
$ws.Range("O3").Value = 0.02703852164627077
$ws.Range("P3").Value = 0.02703852164627077
$ws.Range("S3").Value = 0.02621037747316638
$ws.Range("T3").Value = 0.02621037747316637

# Row 4
$ws.Range("I4").Value = 0.9693716918425304
$ws.Range("J4").Value = 0.9693716918425304
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03070066666666667
$ws.Range("N4").Value = 0.092102
$ws.Range("O4").Value = 0.002835811525282158
$ws.Range("P4").Value = 0.002835811525282158
$ws.Range("Q4").Value = 0.281124756172889
$ws.Range("R4").Value = 2.530122805556
$ws.Range("S4").Value = 0.002748955416009312
$ws.Range("T4").Value = 0.002748955416009312

# Row 5
$ws.Range("G5").Value = 0.2893236666666667
$ws.Range("H5").Value = 0.867971
$ws.Range("I5").Value = 0.03062830815746963
$ws.Range("J5").Value = 0.03062830815746962
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 10.502639
$ws.Range("N5").Value = 31.507917
$ws.Range("O5").Value = 0.9701256668284471
$ws.Range("P5").Value = 0.970125666828447
$ws.Range("Q5").Value = 3.038662025156334
$ws.Range("R5").Value = 27.347958226407
$ws.Range("S5").Value = 0.02971330787509239
$ws.Range("T5").Value = 0.02971330787509238

# Row 6
$ws.Range("G6").Value = 0.2893236666666667
$ws.Range("H6").Value = 0.867971
$ws.Range("I6").Value = 0.03062830815746963
$ws.Range("J6").Value = 0.03062830815746962
$ws.Range("O6").Value = 0.02703852164627077
$ws.Range("P6").Value = 0.02703852164627077
$ws.Range("Q6").Value = 0.08469101658911113
$ws.Range("R6").Value = 0.7622191493020001
$ws.Range("S6").Value = 0.0008281441731043942
$ws.Range("T6").Value = 0.0008281441731043941

# Row 7
$ws.Range("G7").Value = 0.2893236666666667
$ws.Range("H7").Value = 0.867971
$ws.Range("I7").Value = 0.03062830815746963
$ws.Range("J7").Value = 0.03062830815746962
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.03070066666666667
$ws.Range("N7").Value = 0.092102
$ws.Range("O7").Value = 0.002835811525282158
$ws.Range("P7").Value = 0.002835811525282158
$ws.Range("Q7").Value = 0.008882429449111112
$ws.Range("R7").Value = 0.079941865042
$ws.Range("S7").Value = 0.00008685610927284591
$ws.Range("T7").Value = 0.00008685610927284591
